# Update New Orleans shard workbook:
#  1. Insert a new "State" column into hotel_info (between Hotel_Name and City)
#     and populate it with "Louisiana" for the existing hotel row.
#  2. Reorder the worksheet tabs so that review_info comes before hotel_info.

$wb = $excel.ActiveWorkbook

$wsHotel  = $wb.Worksheets.Item("hotel_info")
$wsReview = $wb.Worksheets.Item("review_info")

# --- 1. Insert "State" column (column C) into hotel_info ---
$wsHotel.Columns.Item(3).Insert()
$wsHotel.Cells.Item(1, 3).Value = "State"
$wsHotel.Cells.Item(2, 3).Value = "Louisiana"

# --- 2. Move review_info so it is the first sheet (before hotel_info) ---
$wsReview.Move($wsHotel)
